$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.241.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5069"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09824"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.146"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.85"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.508"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.76%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.877.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.460"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.78%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.14%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06592"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.63"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9985"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.195"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.298.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.28%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.293"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.59%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.575"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.16%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.092.53"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.24%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.30"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.38%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.60"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.85%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.76"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1067"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.72%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.074"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.74%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.653"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.38%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.621"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.543"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +9.08%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06738"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02387"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.90%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2192"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.05%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6395"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.52"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.87%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.59%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.186"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.48%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Frax"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9986"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.78%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.34%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6012"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.28%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.660"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.268"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.08%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.007"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.34%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.28"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.14%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.197"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.38%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06856"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.45%  "
